# nublar_automation.xlsx update
# "21 fertig, 18 weitergearbeitet"
#  - task #21 (Create Volumes, row 20) finished -> informiert/umgesetzt/verifiziert filled in
#  - task #20 (Create Volume Template, row 19) Zustaendig "OlAnt" cleared, Filename filled in
#  - task #18 (Add Hypervisor Cluster Profiles, row 29) Filename filled in
#  - task #17 (Create SP Template, row 31) Filename filled in
#  - active selection moved to K21
#  - hyperlink cell O21 switched from the stray "Hyperlink" cell style to the builtin "Link" style
#  - row 21's explicit row height reset back to the sheet default

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (# 20, "Create Volume Template"): drop the "Zustaendig" (OlAnt) note,
# fill in the Filename column instead.
$ws.Range("H19").Clear()
$ws.Range("L19").Value = "X_20_createvolumetemplate.yml"

# Row 20 (# 21, "Create Volumes"): task is finished -> drop "Zustaendig" (OlAnt) and
# tick off informiert / umgesetzt / verifiziert.
$ws.Range("H20").Clear()
$ws.Range("I20").Value = "X"
$ws.Range("J20").Value = "X"
$ws.Range("K20").Value = "?"

# Row 21: normalise row height back to the sheet default and switch the hyperlink
# cell O21 from the orphan "Hyperlink" cell style to the builtin "Link" style.
$ws.Range("O21").Style = "Link"
$ws.Rows(21).AutoFit()

# That was the only cell still using the custom "Hyperlink" cell style, so it is now
# unused -> drop it (Excel prunes unused cell styles like this automatically).
$wb.Styles("Hyperlink").Delete()

# Row 29 (# 18, "Add Hypervisor Cluster Profiles"): fill in the Filename column.
$ws.Range("L29").Value = "X_18_addhypervisorclusterprofile.yml"

# Row 31 (# 17, "Create SP Template"): fill in the Filename column.
$ws.Range("L31").Value = "X_17_createserverprofiletemplate.yml"

# Move the active selection to K21 (from K20).
$ws.Range("K21").Select()
